# Quarterly indexing bug-fix: shift each date in column A (originally the
# 1st of a quarter-start month: Jan/Apr/Jul/Oct) forward so it lands on the
# 15th of the *following* month instead.
#
# Example: 1988-07-01 (serial 32325) -> 1988-08-15 (serial 32370)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed OLE Automation / Excel-serial epoch (day 0 = 1899-12-30), pinned to
# midnight so arithmetic doesn't pick up the sandbox clock's time-of-day.
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 150 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($serial -eq $null -or $serial -eq "") { continue }

    $d = $epoch.AddDays([double]$serial)
    $d2 = $d.AddMonths(1)
    $newDate = Get-Date -Year $d2.Year -Month $d2.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value2 = $newDate.ToOADate()
}
